$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.952.60'
$ws.Range('E2').Value = '  +0.64%  '
$ws.Range('D3').Value = '2.417.87'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('E4').Value = '  +0.10%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '552.07'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.03%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '137.18'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +0.22%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +0.08%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.589'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +3.38%  '
$ws.Range('E9').Value = '  -1.98%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '5.70'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -2.02%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.147'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -1.08%  '
$ws.Range('B12').Value = 'Cardano'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.354'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -2.00%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '25.51'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +3.26%  '
$ws.Range('D14').Value = '2.848.13'
$ws.Range('E14').Value = '  +0.37%  '
$ws.Range('D15').Value = '59.873.50'
$ws.Range('E15').Value = '  +0.78%  '
$ws.Range('E16').Value = '  -1.31%  '
$ws.Range('D17').Value = '2.433.27'
$ws.Range('E17').Value = '  +1.76%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '11.33'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('E19').Value = '  -0.69%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '329.26'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -1.76%  '
$ws.Range('E21').Value = '  -5.21%  '
$ws.Range('E22').Value = '  +0.10%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '66.66'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +2.85%  '
$ws.Range('E24').Value = '  +1.20%  '
$ws.Range('E25').Value = '  +1.67%  '
$ws.Range('E26').Value = '  +0.12%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '1.37'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +0.46%  '
$ws.Range('D28').Value = '0.0₃0775'
$ws.Range('E28').Value = '  -0.18%  '
$ws.Range('E29').Value = '  -0.92%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '169.31'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -0.73%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '6.11'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -2.46%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '18.64'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -0.28%  '
$ws.Range('E33').Value = '  -1.14%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('E35').Value = '  +0.21%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '4.23'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -1.69%  '
$ws.Range('E37').Value = '  +0.03%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '1.60'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -1.84%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '314.58'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +3.80%  '
$ws.Range('B40').Value = 'PolygonEcosystemToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.408'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -2.74%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '3.67'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -1.87%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '138.70'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -3.01%  '
$ws.Range('E43').Value = '  +0.54%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.0519'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -1.40%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '19.57'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +2.92%  '
$ws.Range('E46').Value = '  +1.37%  '
$ws.Range('E47').Value = '  -0.66%  '
$ws.Range('E48').Value = '  -4.61%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '17.72'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -1.11%  '
$ws.Range('E50').Value = '  +0.25%  '
$ws.Range('E51').Value = '  -0.71%  '
